# Update TPM-derived LR-pair statistics (Sertad1-Ar sheet) with newly
# recomputed values for columns G:T across data rows 2-19, per the new
# TPM run referenced in the commit ("update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 19.010634
$ws.Range("H2").Value = 57.031902
$ws.Range("I2").Value = 0.1669610886139304
$ws.Range("J2").Value = 0.1669610886139304
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2356743333333333
$ws.Range("N2").Value = 0.707023
$ws.Range("O2").Value = 0.04509540876701491
$ws.Range("P2").Value = 0.04509540876701491
$ws.Range("Q2").Value = 4.480318494194
$ws.Range("R2").Value = 40.322866447746
$ws.Range("S2").Value = 0.007529178539230988
$ws.Range("T2").Value = 0.007529178539230989

# Row 3
$ws.Range("G3").Value = 19.010634
$ws.Range("H3").Value = 57.031902
$ws.Range("I3").Value = 0.1669610886139304
$ws.Range("J3").Value = 0.1669610886139304
$ws.Range("O3").Value = 0.790148549139098
$ws.Range("P3").Value = 0.790148549139098
$ws.Range("Q3").Value = 78.5028288835
$ws.Range("R3").Value = 706.5254599515
$ws.Range("S3").Value = 0.1319240619309815
$ws.Range("T3").Value = 0.1319240619309815

# Row 4
$ws.Range("G4").Value = 19.010634
$ws.Range("H4").Value = 57.031902
$ws.Range("I4").Value = 0.1669610886139304
$ws.Range("J4").Value = 0.1669610886139304
$ws.Range("M4").Value = 0.861036
$ws.Range("N4").Value = 2.583108
$ws.Range("O4").Value = 0.1647560420938871
$ws.Range("P4").Value = 0.1647560420938871
$ws.Range("Q4").Value = 16.368840256824
$ws.Range("R4").Value = 147.319562311416
$ws.Range("S4").Value = 0.02750784814371793
$ws.Range("T4").Value = 0.02750784814371793

# Row 5
$ws.Range("I5").Value = 0.07884104858418527
$ws.Range("J5").Value = 0.07884104858418528
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2356743333333333
$ws.Range("N5").Value = 0.707023
$ws.Range("O5").Value = 0.04509540876701491
$ws.Range("P5").Value = 0.04509540876701491
$ws.Range("Q5").Value = 2.115660666840555
$ws.Range("R5").Value = 19.040946001565
$ws.Range("S5").Value = 0.003555369313523917
$ws.Range("T5").Value = 0.003555369313523918

# Row 6
$ws.Range("I6").Value = 0.07884104858418527
$ws.Range("J6").Value = 0.07884104858418528
$ws.Range("O6").Value = 0.790148549139098
$ws.Range("P6").Value = 0.790148549139098
$ws.Range("S6").Value = 0.06229614015139913
$ws.Range("T6").Value = 0.06229614015139914

# Row 7
$ws.Range("I7").Value = 0.07884104858418527
$ws.Range("J7").Value = 0.07884104858418528
$ws.Range("M7").Value = 0.861036
$ws.Range("N7").Value = 2.583108
$ws.Range("O7").Value = 0.1647560420938871
$ws.Range("P7").Value = 0.1647560420938871
$ws.Range("Q7").Value = 7.72956465886
$ws.Range("R7").Value = 69.56608192974001
$ws.Range("S7").Value = 0.01298953911926223
$ws.Range("T7").Value = 0.01298953911926223

# Row 8
$ws.Range("G8").Value = 21.69976666666667
$ws.Range("H8").Value = 65.0993
$ws.Range("I8").Value = 0.1905784239144757
$ws.Range("J8").Value = 0.1905784239144757
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2356743333333333
$ws.Range("N8").Value = 0.707023
$ws.Range("O8").Value = 0.04509540876701491
$ws.Range("P8").Value = 0.04509540876701491
$ws.Range("Q8").Value = 5.114078042655555
$ws.Range("R8").Value = 46.02670238389999
$ws.Range("S8").Value = 0.008594211928596732
$ws.Range("T8").Value = 0.008594211928596733

# Row 9
$ws.Range("G9").Value = 21.69976666666667
$ws.Range("H9").Value = 65.0993
$ws.Range("I9").Value = 0.1905784239144757
$ws.Range("J9").Value = 0.1905784239144757
$ws.Range("O9").Value = 0.790148549139098
$ws.Range("P9").Value = 0.790148549139098
$ws.Range("Q9").Value = 89.60737813611111
$ws.Range("R9").Value = 806.4664032249999
$ws.Range("S9").Value = 0.150585265153239
$ws.Range("T9").Value = 0.150585265153239

# Row 10
$ws.Range("G10").Value = 21.69976666666667
$ws.Range("H10").Value = 65.0993
$ws.Range("I10").Value = 0.1905784239144757
$ws.Range("J10").Value = 0.1905784239144757
$ws.Range("M10").Value = 0.861036
$ws.Range("N10").Value = 2.583108
$ws.Range("O10").Value = 0.1647560420938871
$ws.Range("P10").Value = 0.1647560420938871
$ws.Range("Q10").Value = 18.6842802916
$ws.Range("R10").Value = 168.1585226244
$ws.Range("S10").Value = 0.03139894683264002
$ws.Range("T10").Value = 0.03139894683264003

# Row 11
$ws.Range("G11").Value = 6.543946666666667
$ws.Range("H11").Value = 19.63184
$ws.Range("I11").Value = 0.05747227889917651
$ws.Range("J11").Value = 0.05747227889917651
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.2356743333333333
$ws.Range("N11").Value = 0.707023
$ws.Range("O11").Value = 0.04509540876701491
$ws.Range("P11").Value = 0.04509540876701491
$ws.Range("Q11").Value = 1.542240268035556
$ws.Range("R11").Value = 13.88016241232
$ws.Range("S11").Value = 0.00259173590973025
$ws.Range("T11").Value = 0.00259173590973025

# Row 12
$ws.Range("G12").Value = 6.543946666666667
$ws.Range("H12").Value = 19.63184
$ws.Range("I12").Value = 0.05747227889917651
$ws.Range("J12").Value = 0.05747227889917651
$ws.Range("O12").Value = 0.790148549139098
$ws.Range("P12").Value = 0.790148549139098
$ws.Range("Q12").Value = 27.02268243111111
$ws.Range("R12").Value = 243.20414188
$ws.Range("S12").Value = 0.04541163778790192
$ws.Range("T12").Value = 0.04541163778790192

# Row 13
$ws.Range("G13").Value = 6.543946666666667
$ws.Range("H13").Value = 19.63184
$ws.Range("I13").Value = 0.05747227889917651
$ws.Range("J13").Value = 0.05747227889917651
$ws.Range("M13").Value = 0.861036
$ws.Range("N13").Value = 2.583108
$ws.Range("O13").Value = 0.1647560420938871
$ws.Range("P13").Value = 0.1647560420938871
$ws.Range("Q13").Value = 5.63457366208
$ws.Range("R13").Value = 50.71116295872
$ws.Range("S13").Value = 0.009468905201544345
$ws.Range("T13").Value = 0.009468905201544347

# Row 14
$ws.Range("G14").Value = 45.17409633333333
$ws.Range("H14").Value = 135.522289
$ws.Range("I14").Value = 0.3967419656263906
$ws.Range("J14").Value = 0.3967419656263906
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.2356743333333333
$ws.Range("N14").Value = 0.707023
$ws.Range("O14").Value = 0.04509540876701491
$ws.Range("P14").Value = 0.04509540876701491
$ws.Range("Q14").Value = 10.64637503729411
$ws.Range("R14").Value = 95.817375335647
$ws.Range("S14").Value = 0.01789124111495106
$ws.Range("T14").Value = 0.01789124111495106

# Row 15
$ws.Range("G15").Value = 45.17409633333333
$ws.Range("H15").Value = 135.522289
$ws.Range("I15").Value = 0.3967419656263906
$ws.Range("J15").Value = 0.3967419656263906
$ws.Range("O15").Value = 0.790148549139098
$ws.Range("P15").Value = 0.790148549139098
$ws.Range("Q15").Value = 186.5426663004722
$ws.Range("R15").Value = 1678.88399670425
$ws.Range("S15").Value = 0.3134850885222864
$ws.Range("T15").Value = 0.3134850885222865

# Row 16
$ws.Range("G16").Value = 45.17409633333333
$ws.Range("H16").Value = 135.522289
$ws.Range("I16").Value = 0.3967419656263906
$ws.Range("J16").Value = 0.3967419656263906
$ws.Range("M16").Value = 0.861036
$ws.Range("N16").Value = 2.583108
$ws.Range("O16").Value = 0.1647560420938871
$ws.Range("P16").Value = 0.1647560420938871
$ws.Range("Q16").Value = 38.896523210468
$ws.Range("R16").Value = 350.068708894212
$ws.Range("S16").Value = 0.06536563598915313
$ws.Range("T16").Value = 0.06536563598915314

# Row 17
$ws.Range("G17").Value = 12.45716666666667
$ws.Range("H17").Value = 37.3715
$ws.Range("I17").Value = 0.1094051943618415
$ws.Range("J17").Value = 0.1094051943618415
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.2356743333333333
$ws.Range("N17").Value = 0.707023
$ws.Range("O17").Value = 0.04509540876701491
$ws.Range("P17").Value = 0.04509540876701491
$ws.Range("Q17").Value = 2.935834449388889
$ws.Range("R17").Value = 26.4225100445
$ws.Range("S17").Value = 0.004933671960981958
$ws.Range("T17").Value = 0.004933671960981959

# Row 18
$ws.Range("G18").Value = 12.45716666666667
$ws.Range("H18").Value = 37.3715
$ws.Range("I18").Value = 0.1094051943618415
$ws.Range("J18").Value = 0.1094051943618415
$ws.Range("O18").Value = 0.790148549139098
$ws.Range("P18").Value = 0.790148549139098
$ws.Range("Q18").Value = 51.44083165277777
$ws.Range("R18").Value = 462.967484875
$ws.Range("S18").Value = 0.0864463555932901
$ws.Range("T18").Value = 0.0864463555932901

# Row 19
$ws.Range("G19").Value = 12.45716666666667
$ws.Range("H19").Value = 37.3715
$ws.Range("I19").Value = 0.1094051943618415
$ws.Range("J19").Value = 0.1094051943618415
$ws.Range("M19").Value = 0.861036
$ws.Range("N19").Value = 2.583108
$ws.Range("O19").Value = 0.1647560420938871
$ws.Range("P19").Value = 0.1647560420938871
$ws.Range("Q19").Value = 10.726068958
$ws.Range("R19").Value = 96.53462062200001
$ws.Range("S19").Value = 0.01802516680756946
$ws.Range("T19").Value = 0.01802516680756946

Write-Output "Applied 208 cell updates to Sheet1 (rows 2-19, cols G:T)"
